$wb = $excel.ActiveWorkbook

# --- Sheet "All Orders": insert a new order row at the top of the data (row 2) ---
$ws = $wb.Worksheets.Item("All Orders")
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "2026-01-13 22:38"
$ws.Cells.Item(2, 3).Value = "Phantom"
# D2 ("420") and J2 ("2026-01-14") look numeric/date-like, so force text
# formatting before assignment to keep them as literal strings, matching
# every other Flat No / Collection Date cell in the sheet.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "420"
# E2 (Phone) is blank - leave the cell unset so it stays empty, as with
# the other blank cells produced by the row-insert shift.
$ws.Cells.Item(2, 6).Value = "Upma x1"
$ws.Cells.Item(2, 7).Value = 30
$ws.Cells.Item(2, 8).Value = "NEW"
$ws.Cells.Item(2, 9).Value = "PENDING"
$ws.Cells.Item(2, 10).NumberFormat = "@"
$ws.Cells.Item(2, 10).Value = "2026-01-14"
$ws.Cells.Item(2, 11).Value = "15:38"
# L2, M2, N2 (Notes / Cancel Reason / Feedback) are blank - leave unset.

# --- Sheet "Daily Summary": refresh the daily roll-up totals ---
$ws2 = $wb.Worksheets.Item("Daily Summary")
$ws2.Cells.Item(2, 2).Value = 11
$ws2.Cells.Item(2, 5).Value = 275
$ws2.Cells.Item(2, 7).Value = 275
